$d = $word.ActiveDocument

# 1. SCALLOP-Seq merge (gramStart/gramEnd removal, no visible text change)
$d.Content.Find.Execute("the SCALLOP-Seq(", $true, $false, $false, $false, $false, $true, 1, $false, "the SCALLOP-Seq(", 2)

# 2. URL change
$d.Content.Find.Execute("https://cambridge-ceu.github.io/csd3/", $true, $false, $false, $false, $false, $true, 1, $false, "https://cambridge-ceu.github.io/csd3/systems/ceuadmin.html", 2)
